# Remove needless imports on org.dozer
#
# Sheet1 ("Environment") lists the imports used by the mapping. The
# "org.dozer" import (row 10, column D) is no longer needed, so delete
# the entire row. Excel automatically shifts the rows below it up by
# one, which also fixes up the dependent merged-cell ranges
# (C8:C10 -> C8:C9, C13:H13 -> C12:H12), the sheet dimension, and the
# now-unused "org.dozer" shared string.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(10).Delete() | Out-Null

# Switch focus to Sheet1 (it becomes the active tab instead of Sheet2)
# and land the selection on D10, which after the shift holds
# "org.openl.rules.mapping.to".
$ws1.Activate() | Out-Null
$ws1.Range("D10").Select() | Out-Null
